# Fixed a bug in respin
# Re-orders the reel-weight rows (rows 2-21, columns A:F) on the active sheet
# so that each row's data lands on its corrected row position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(101,9,30,15,60,15),
    @(1001,18,30,75,60,72),
    @(501,9,52,30,75,45),
    @(902,1,0,0,0,0),
    @(301,6,45,30,60,45),
    @(401,9,48,67,75,45),
    @(1203,3,15,15,15,15),
    @(701,3,90,45,97,15),
    @(601,9,60,67,60,42),
    @(201,9,30,15,45,30),
    @(1201,2,10,10,10,10),
    @(901,16,15,45,60,60),
    @(801,3,67,65,52,45),
    @(1202,2,10,10,10,10),
    @(502,0,4,0,0,0),
    @(802,0,4,5,4,0),
    @(2,0,2,2,2,2),
    @(3,0,3,3,3,3),
    @(1101,0,15,30,30,0),
    @(1,0,2,2,2,2)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Count; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $values[$j]
    }
}
